# Update column F (dSF) values to match repulled data / recalculated means
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -14
$ws.Range("F3").Value = -6
$ws.Range("F4").Value = -4
$ws.Range("F6").Value = -10
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -4
$ws.Range("F9").Value = -2
$ws.Range("F10").Value = -8
$ws.Range("F13").Value = -5
$ws.Range("F15").Value = 1
$ws.Range("F18").Value = 1
